$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.683.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.735.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9983"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9963"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3732"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3406"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.189"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07468"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9958"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.444"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.091"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.731.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.98%  "

$ws.Range("E17").Value = "  +2.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06723"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9969"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.234"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "26.664.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.430"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.507"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +27.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.435"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.919.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.051"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08595"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.698"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.416"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02351"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2176"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06234"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.512"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.223"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6247"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9967"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.908"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6076"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.054"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07198"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.33%  "
